# Applies two changes described by the commit diff:
#  1) Merge the two "MON Nov 06" / " 11:55:05 PST 2017" runs into a single run.
#  2) Append a new transaction block (MON Dec 04 09:38:08 PST 2017 / MAMATHA / CHICK IN)
#     right after the last "Amount Received mode ... - CASH" paragraph.

$d = $word.ActiveDocument

# --- Change 1: merge the split date/time runs into one run -----------------
$d.Content.Find.Execute("MON Nov 06 11:55:05 PST 2017", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "MON Nov 06 11:55:05 PST 2017", 2) | Out-Null

# --- Change 2: append the new purchase-detail block -------------------------
# Find the last paragraph that holds "Amount Received mode ... - CASH" (the
# most recent transaction at the end of the document) and add the new block
# right after it.
$paras = $d.Paragraphs
$targetIndex = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    $txt = $paras.Item($i).Range.Text
    if ($txt -like "Amount Received mode*- CASH*") {
        $targetIndex = $i
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the 'Amount Received mode ... - CASH' paragraph."
}

$r = $paras.Item($targetIndex).Range
$r.Collapse(0)

# Blank separator paragraph.
$r.InsertAfter("`r")
$r.Collapse(0)

# Date/time line (kept as plain text; the source author's original had this
# split across two runs with identical formatting, which collapse together).
$r.InsertAfter("MON Dec 04")
$r.Collapse(0)
$r.InsertAfter(" 09:38:08 PST 2017")
$r.Collapse(0)
$r.InsertAfter("`r")
$r.Collapse(0)

# Person Name line.
$r.InsertAfter("Person Name`t`t`t`t- P")
$r.Collapse(0)
$r.InsertAfter("`r")
$r.Collapse(0)

# Bill number line.
$r.InsertAfter("Bill number`t`t`t`t- 1776")
$r.Collapse(0)
$r.InsertAfter("`r")
$r.Collapse(0)

# Separator line.
$r.InsertAfter("---------------------------------------------------------------")
$r.Collapse(0)
$r.InsertAfter("`r")
$r.Collapse(0)

# Item Name line.
$r.InsertAfter("Item Name`t`t`t`t- CARROT")
$r.Collapse(0)
$r.InsertAfter("`r")
$r.Collapse(0)

# Number of Pockets line.
$r.InsertAfter("Number of Pockets`t`t`t- 1")
$r.Collapse(0)
$r.InsertAfter("`r")
$r.Collapse(0)

# Number of KGs line.
$r.InsertAfter("Number of KGs`t`t`t- 41")
$r.Collapse(0)
$r.InsertAfter("`r")
$r.Collapse(0)

# Rate line.
$r.InsertAfter("Rate`t`t`t`t`t- 58")
$r.Collapse(0)
$r.InsertAfter("`r")
$r.Collapse(0)

# Total Price line.
$r.InsertAfter("Total Price`t`t`t`t- 2378.0")
$r.Collapse(0)
$r.InsertAfter("`r")
$r.Collapse(0)

# Amount balance line (bold).
$r.InsertAfter("Amount balance`t`t`t- 17434.0")
$r.Collapse(0)
$boldRange = $d.Range($r.Start - ("Amount balance`t`t`t- 17434.0").Length, $r.Start)
$boldRange.Font.Bold = 1
$r.InsertAfter("`r")
$r.Collapse(0)

# Trailing blank paragraph.
$r.InsertAfter("`r")

Write-Host "Edit applied."
